$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell -> new value (derived from the authoritative XML diff).
$updates = @{
    "A1" = 45311
    "D28" = 469.482
    "D29" = 536.554
    "D30" = 613.997
    "D31" = 718.271
    "D32" = 785.847
    "D33" = 936.438
    "D34" = 959.213
    "D35" = 1070.066
    "D36" = 1295.822
    "D37" = 1462.864
    "D38" = 1708.361
    "D39" = 1885.526
    "D40" = 2050.034
    "D41" = 2303.124
    "D42" = 2457.509
    "D43" = 2680.229
    "D49" = 1505.89
    "D50" = 1513.482
    "D51" = 1857.684
    "D52" = 1986.763
    "D53" = 2116.086
    "D54" = 2381.58
    "D55" = 2771.34
    "D56" = 3102.889
    "D57" = 3548.327
    "D58" = 3998.828
    "D59" = 4340.497
    "D60" = 4884.644
    "D61" = 5388.291
    "D62" = 5707.184
    "D68" = 607.671
    "D69" = 625.134
    "D70" = 707.388
    "D71" = 824.064
    "D72" = 100
    "D73" = 1104.993
    "D74" = 1237.612
    "D75" = 1311.008
    "D76" = 1543.85
    "D77" = 1789.351
    "D78" = 1993.592
    "D79" = 2234.28
    "D80" = 2492.942
    "D81" = 2640.49
    "D82" = 2862.453
    "D83" = 3125.665
    "D89" = 259.878
    "D90" = 312.062
    "D91" = 337.37
    "D92" = 418.61
    "D93" = 456.827
    "D94" = 492.008
    "D95" = 539.842
    "D96" = 688.913
    "D97" = 749.146
    "D98" = 878.223
    "D99" = 992.114
    "D100" = 1068.042
    "D106" = 2406.888
    "D107" = 2558.744
    "D108" = 2708.068
    "D109" = 2961.16
    "D110" = 2978.875
    "D111" = 3510.364
    "D112" = 4155.742
    "D113" = 4606.245
    "D114" = 5036.501
    "D115" = 5560.397
    "D116" = 6226.021
    "D117" = 6577.82
    "D118" = 7339.619
    "D119" = 7782.528
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Output "Updated $($updates.Count) cells."
